$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.459.99"
$ws.Range("E2").Value = "  +3.40%  "
$ws.Range("D3").Value = "1.602.86"
$ws.Range("E3").Value = "  +3.17%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'212.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'0.518"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.09%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'26.88"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.72%  "
$ws.Range("D9").Value = "'43.43"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("E10").Value = "  +2.71%  "
$ws.Range("D11").Value = "'0.0598"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.74%  "
$ws.Range("D12").Value = "'0.0910"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").Value = "1.831.06"
$ws.Range("E13").Value = "  +3.05%  "
$ws.Range("D14").Value = "1.616.34"
$ws.Range("E14").Value = "  +4.14%  "
$ws.Range("D15").Value = "29.527.04"
$ws.Range("E15").Value = "  +3.58%  "
$ws.Range("E16").Value = "  +4.90%  "
$ws.Range("D17").Value = "'3.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.10%  "
$ws.Range("D18").Value = "'63.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.83%  "
$ws.Range("D19").Value = "'239.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.59%  "
$ws.Range("D20").Value = "'7.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.49%  "
$ws.Range("D21").Value = "0.0₃0693"
$ws.Range("E21").Value = "  +3.30%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "'4.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.44%  "
$ws.Range("D24").Value = "'9.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.09%  "
$ws.Range("D26").Value = "'154.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.20%  "
$ws.Range("D27").Value = "'15.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.73%  "
$ws.Range("E28").Value = "  +5.33%  "
$ws.Range("D29").Value = "'6.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.25%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("E31").Value = "  +2.66%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("E33").Value = "  +2.37%  "
$ws.Range("D34").Value = "1.424.39"
$ws.Range("E34").Value = "  +2.65%  "
$ws.Range("E35").Value = "  +3.51%  "
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("E37").Value = "  +1.95%  "
$ws.Range("D38").Value = "'2.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.16%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  +2.06%  "
$ws.Range("E41").Value = "  +3.48%  "
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("D43").Value = "'53.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +21.43%  "
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").Value = "'0.792"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.12%  "
$ws.Range("E46").Value = "  +1.65%  "
$ws.Range("D47").Value = "'65.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.69%  "
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("D49").Value = "1.743.04"
$ws.Range("E49").Value = "  +2.64%  "
$ws.Range("D50").Value = "'86.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.50%  "
$ws.Range("D51").Value = "'0.836"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.55%  "
